$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openTickets")

# Fix the malformed JSON (missing colon before "Service Anfrage") in the
# ticketDescriptionHighlighting value held in G3.
$g3Old = $ws.Range("G3").Value()
$ws.Range("G3").Value = $g3Old.Replace('"key" "Service Anfrage"', '"key":"Service Anfrage"')

# Update the active selection on the sheet to G5 (was I5).
$ws.Range("G5").Select()
